# Add a "Transcriptomics" tag to the minSCe single cell sequencing template.
# This mirrors the existing "single-cell RNA sequencing assay" tag (columns
# B/C of the TAGS block on the isa_template sheet) by adding a third tag in
# column E, rows 13-17:
#   E13 Tags                       -> Transcriptomics
#   E14 Tags Term Accession Number -> https://bioregistry.io/NCIT:C153189
#   E15 Tags Term Source REF       -> NCIT
#   E16 Comment[description]       -> "A study of the complete set of RNA
#                                       transcripts that are produced by the
#                                       genome, under specific circumstances
#                                       or in a specific cell." []
#   E17 Comment[isObsolete]        -> false

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("isa_template")

$ws.Range("E13").Value = "Transcriptomics"
$ws.Range("E14").Value = "https://bioregistry.io/NCIT:C153189"
$ws.Range("E15").Value = "NCIT"
$ws.Range("E16").Value = '"A study of the complete set of RNA transcripts that are produced by the genome, under specific circumstances or in a specific cell." []'
$ws.Range("E17").Value = "'false"

$ws.Range("E13:E17").Select()
